$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Participant p3's measurements (row 4 on the sheet; row 1 is the header,
# so row 4 = IDN 3) were recorded only for the second ascent via the IMU,
# so that participant's row is dropped from the demographic table and all
# later rows shift up to close the gap. The running IDN numbers in column
# A are left exactly as they were (they already form an unbroken 1..N
# sequence), so only columns B:H need to be pulled up from the row below.
for ($r = 4; $r -le 17; $r++) {
    for ($c = 2; $c -le 8; $c++) {
        $below = $ws.Cells.Item($r + 1, $c).Value()
        $ws.Cells.Item($r, $c).Value = $below
    }
}

# The old last row (18) is now a duplicate of row 17; clear it so it drops
# out of the sheet's used range entirely.
$ws.Range("A18:H18").ClearContents()

$ws.Range("J14").Select() | Out-Null
